$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "Word Scaling Analysis"
$ws.Range("F2").Value = 4

$ws.Range("G2").Value = "Long English Scaling"
$ws.Range("G3").Value = "Nonsense Scaling"
$ws.Range("G4").Value = "English Scaling"
$ws.Range("G5").Value = "Hebrew Scaling"

$ws.Range("H4").Value = 3
$ws.Range("H5").Value = 4

$ws.Range("I2").Value = "red"
$ws.Range("I3").Value = "cyan"
$ws.Range("I4").Value = "blue"
$ws.Range("I5").Value = "black"
